$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new values
$ws.Range("A2").Value = "PROD002"
$ws.Range("B2").Value = "SICA"
$ws.Range("C2").Value = "PRODUCTO01"
$ws.Range("F2").Value = "ferreteria_general"

# Delete rows 3 to 5 entirely (shifts cells up, removing them from the sheet)
$ws.Range("A3:F5").Delete()
